$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7896
$ws.Range("F5").Value = 7896
$ws.Range("F8").Value = 2139
$ws.Range("F9").Value = 8653
$ws.Range("F12").Value = 94
$ws.Range("F13").Value = 5789
$ws.Range("F15").Value = 2784
$ws.Range("F16").Value = 1208
$ws.Range("F17").Value = 420
$ws.Range("F20").Value = 624
$ws.Range("F21").Value = 98
$ws.Range("F22").Value = 3966
$ws.Range("F29").Value = 5645
$ws.Range("F34").Value = 409
$ws.Range("F35").Value = 3019
$ws.Range("F36").Value = 1544
$ws.Range("F38").Value = 1424
$ws.Range("F39").Value = 5662
$ws.Range("F43").Value = 3651
$ws.Range("F44").Value = 23
$ws.Range("F46").Value = 2350
$ws.Range("F50").Value = 25
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 19
$ws.Range("F9").Value = 33
$ws.Range("F10").Value = 133
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 7896
$ws.Range("F6").Value = 7896
$ws.Range("F9").Value = 2139
$ws.Range("F10").Value = 8653
$ws.Range("F13").Value = 94
$ws.Range("F14").Value = 5789
$ws.Range("F16").Value = 2784
$ws.Range("F17").Value = 1208
$ws.Range("F18").Value = 420
$ws.Range("F22").Value = 624
$ws.Range("F23").Value = 98
$ws.Range("F24").Value = 3966
$ws.Range("F30").Value = 5645
$ws.Range("F34").Value = 409
$ws.Range("F35").Value = 3019
$ws.Range("F36").Value = 1544
$ws.Range("F39").Value = 1424
$ws.Range("F41").Value = 5662
$ws.Range("F44").Value = 3651
$ws.Range("F46").Value = 33
$ws.Range("F47").Value = 2350
